$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.347.17'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '1.870.69'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.14'
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4673'
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2844'
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06554'
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.89'
$ws.Range("E10").Value = '  +7.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07906'
$ws.Range("E11").Value = '  +1.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.32'
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D13").Value = '1.879.26'
$ws.Range("E13").Value = '  -0.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.155'
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6778'
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '285.89'
$ws.Range("E16").Value = '  +1.67%  '
$ws.Range("D17").Value = '30.376.26'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.550'
$ws.Range("E18").Value = '  +3.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9986'
$ws.Range("E19").Value = '  -0.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.71'
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("D21").Value = '2.113.75'
$ws.Range("E21").Value = '  -1.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007283'
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  -0.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.199'
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.323'
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.23'
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.15'
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.917'
$ws.Range("E28").Value = '  -3.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.359'
$ws.Range("E29").Value = '  -2.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09698'
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.427'
$ws.Range("E31").Value = '  -0.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.473'
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.111'
$ws.Range("E33").Value = '  -2.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04721'
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.123'
$ws.Range("E35").Value = '  +2.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7044'
$ws.Range("E36").Value = '  -0.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.717'
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01868'
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.424'
$ws.Range("E39").Value = '  -3.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.543'
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.90'
$ws.Range("E41").Value = '  +1.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.947'
$ws.Range("E42").Value = '  -0.78%  '
$ws.Range("E43").Value = '  -1.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4189'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.31'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.213'
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.306'
$ws.Range("E48").Value = '  +2.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '930.22'
$ws.Range("E49").Value = '  -5.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.15'
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1137'
$ws.Range("E51").Value = '  -2.70%  '
